# Edit readme.docx per the commit diff:
#  1. Move the "_GoBack" bookmark from the end of the "English (United
#     States) spelling" paragraph to the end of the "Type './test'"
#     paragraph (last paragraph of the Terminal how-to section).
#  2. Expand the "g++ ..." compile-command paragraph: add a
#     "calculator.h calculator.cpp" pair to the command line and wrap
#     the non-dictionary words (entity.h, fileImporter.h, calculator.h,
#     std, c++) in spell-check-exception (proofErr) markers.
#  3. Mark "AoE" and "Korei" in the page header as spell-check
#     exceptions (proofErr spellStart/spellEnd) by splitting the runs
#     that contain them.

$d = $word.ActiveDocument

# --- 1a. Strip the bookmark off the "English (United States) spelling." paragraph.
$pEnglish = $d.Paragraphs(7)
if ($pEnglish.Range.Text -notmatch "English") {
    throw "Paragraph 7 is not the expected 'English (United States) spelling' paragraph"
}
$xmlEnglish = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>In regards to the name, use English (United States) spelling.</w:t></w:r></w:p>'
$pEnglish.Range.InsertXML($xmlEnglish)

# --- 2. Rebuild the "g++ ..." NumberedList paragraph with the extra
#        calculator.h/calculator.cpp tokens and proofErr wraps.
$pGxx = $d.Paragraphs(13)
if ($pGxx.Range.Text -notmatch "g\+\+ main.cpp") {
    throw "Paragraph 13 is not the expected 'g++ main.cpp ...' paragraph"
}
$xmlGxx = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NumberedList"/><w:rPr><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">Type “g++ main.cpp </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>entity.h</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> entity.cpp </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>fileImporter.h</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> fileImporter.cpp </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>calculator.h</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> calculator.cpp </w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">–o test – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>std</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>c++</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>11”</w:t></w:r></w:p>'
$pGxx.Range.InsertXML($xmlGxx)

# --- 1b. Re-add the bookmark at the end of the "Type './test'" paragraph.
$pTest = $d.Paragraphs(14)
if ($pTest.Range.Text -notmatch "\./test") {
    throw "Paragraph 14 is not the expected 'Type ./test' paragraph"
}
$xmlTest = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NumberedList"/><w:rPr><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>Type “./test”</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pTest.Range.InsertXML($xmlTest)

# --- 3. Header: mark "AoE" and "Korei" as proofing exceptions.
$hdr = $d.Sections(1).Headers(1)
$pHeader = $hdr.Range.Paragraphs(1)
if ($pHeader.Range.Text -notmatch "Korei Khan") {
    throw "Header paragraph is not the expected title/author line"
}
$xmlHeader = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>AoE</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> II: The Bo</w:t></w:r><w:r><w:t>ard Game -</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Combat results calculator</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">Written by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Korei</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Khan</w:t></w:r></w:p>'
$pHeader.Range.InsertXML($xmlHeader)

Write-Output "done"
